$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.030.77'
$ws.Range('E2').Value = '  -4.33%  '
$ws.Range('D3').Value = '3.280.64'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '558.07'
$ws.Range('E5').Value = '  -3.03%  '
$ws.Range('D6').Value = '185.47'
$ws.Range('E6').Value = '  -2.25%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -2.35%  '
$ws.Range('D9').Value = '3.274.45'
$ws.Range('E9').Value = '  -4.79%  '
$ws.Range('D10').Value = '0.188'
$ws.Range('E10').Value = '  -7.58%  '
$ws.Range('E11').Value = '  -4.47%  '
$ws.Range('D12').Value = '47.44'
$ws.Range('E12').Value = '  -7.13%  '
$ws.Range('E13').Value = '  -5.63%  '
$ws.Range('D14').Value = '8.60'
$ws.Range('E14').Value = '  -4.98%  '
$ws.Range('D15').Value = '632.55'
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('D16').Value = '3.806.74'
$ws.Range('E16').Value = '  -4.65%  '
$ws.Range('D17').Value = '65.924.94'
$ws.Range('E17').Value = '  -4.26%  '
$ws.Range('E18').Value = '  -0.63%  '
$ws.Range('E19').Value = '  -3.19%  '
$ws.Range('D20').Value = '3.278.32'
$ws.Range('E20').Value = '  -4.59%  '
$ws.Range('D21').Value = '11.38'
$ws.Range('E21').Value = '  -6.52%  '
$ws.Range('D22').Value = '0.906'
$ws.Range('E22').Value = '  -3.24%  '
$ws.Range('D23').Value = '17.87'
$ws.Range('E23').Value = '  +1.27%  '
$ws.Range('D24').Value = '106.72'
$ws.Range('E24').Value = '  +8.69%  '
$ws.Range('D25').Value = '4.93'
$ws.Range('E25').Value = '  -7.04%  '
$ws.Range('E26').Value = '  -6.36%  '
$ws.Range('E27').Value = '  -5.75%  '
$ws.Range('D28').Value = '9.55'
$ws.Range('E28').Value = '  -2.65%  '
$ws.Range('D29').Value = '8.71'
$ws.Range('E29').Value = '  -5.02%  '
$ws.Range('D30').Value = '30.45'
$ws.Range('E30').Value = '  -5.16%  '
$ws.Range('D31').Value = '4.07'
$ws.Range('E31').Value = '  -5.37%  '
$ws.Range('D32').Value = '6.27'
$ws.Range('E32').Value = '  -5.73%  '
$ws.Range('D33').Value = '11.04'
$ws.Range('E33').Value = '  -4.09%  '
$ws.Range('D35').Value = '540.50'
$ws.Range('E35').Value = '  +10.23%  '
$ws.Range('D36').Value = '57.31'
$ws.Range('E36').Value = '  -5.94%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').Value = '3.666.30'
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('E39').Value = '  -0.27%  '
$ws.Range('D40').Value = '0.0₃0735'
$ws.Range('E40').Value = '  -6.28%  '
$ws.Range('E41').Value = '  -0.83%  '
$ws.Range('D42').Value = '2.74'
$ws.Range('E42').Value = '  -5.61%  '
$ws.Range('E43').Value = '  -3.89%  '
$ws.Range('D44').Value = '32.59'
$ws.Range('E44').Value = '  -4.23%  '
$ws.Range('E45').Value = '  -7.91%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '3.27'
$ws.Range('E46').Value = '  -2.12%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '0.0414'
$ws.Range('E47').Value = '  -4.58%  '
$ws.Range('E48').Value = '  -6.12%  '
$ws.Range('E49').Value = '  -3.09%  '
$ws.Range('D50').Value = '0.998'
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('E51').Value = '  +2.03%  '
